$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 29 (Giovanni, 21/01/2018 activity) needs to become row 30.
# Two brand new rows are inserted: a new row 29 (Mirko, 20/01/2018) and a
# new row 31 (Mirko, 21/01/2018) right after the (shifted) old row.

# Insert a fresh row at 29 - this pushes the existing row 29 down to row 30.
$ws.Rows.Item(29).Insert()

# Insert another fresh row at 31 (right after the old data, now at row 30).
$ws.Rows.Item(31).Insert()

# --- New row 29: Mirko, 20/01/2018, "progresso e test/studio" ---
$ws.Range("A29").Value = 43120
$ws.Range("B29").Value = "Mirko"
$ws.Range("C29").Value = "progresso e test/studio"
$ws.Range("D29").Value = 0.08333333333333333
$ws.Range("E29").Value = 0.20833333333333334

# --- New row 31: Mirko, 21/01/2018, "sistemato metodi clear back display" ---
$ws.Range("A31").Value = 43121
$ws.Range("B31").Value = "Mirko"
$ws.Range("C31").Value = "sistemato metodi clear back display"
$ws.Range("D31").Value = 0.0625
$ws.Range("E31").Value = 0.0625

# Row 31 wraps onto two lines of text, same auto height as similar rows.
$ws.Rows.Item(31).RowHeight = 28.8

# Update the on-screen selection to match where the user ended up editing.
$ws.Range("C32").Select() | Out-Null
